$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update column B (Longitude longK)
$ws.Range("B3").Value = 50
$ws.Range("B4").Value = 39
$ws.Range("B5").Value = 57
$ws.Range("B6").Value = 55
$ws.Range("B7").Value = 35
$ws.Range("B8").Value = 35
$ws.Range("B9").Value = 44
$ws.Range("B10").Value = 59

# Update column C (Latitude latK)
$ws.Range("C2").Value = 22
$ws.Range("C3").Value = 22
$ws.Range("C4").Value = 22
$ws.Range("C6").Value = 27
$ws.Range("C9").Value = 30
$ws.Range("C10").Value = 30

# Update the active cell selection to B4
$ws.Range("B4").Select()
